$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price/Volume columns remain text (matches original inline-string cells)
# so values like "1.00" or "27.90" keep their trailing zeros instead of
# being auto-converted to numbers by Excel.
$ws.Range("D2:E51").NumberFormat = "@"

# --- Cell value updates from diff ---
$ws.Range("D2").Value = "65.259.10"
$ws.Range("E2").Value = "  +2.32%  "
$ws.Range("D3").Value = "3.167.22"
$ws.Range("E3").Value = "  +3.81%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "579.13"
$ws.Range("E5").Value = "  +4.08%  "
$ws.Range("D6").Value = "151.31"
$ws.Range("E6").Value = "  +6.74%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("D8").Value = "3.164.20"
$ws.Range("E8").Value = "  +3.73%  "
$ws.Range("D9").Value = "0.531"
$ws.Range("E9").Value = "  +2.19%  "
$ws.Range("E10").Value = "  +6.22%  "
$ws.Range("D11").Value = "6.21"
$ws.Range("E11").Value = "  +0.40%  "
$ws.Range("D12").Value = "0.503"
$ws.Range("E12").Value = "  +5.35%  "
$ws.Range("E13").Value = "  +17.17%  "
$ws.Range("D14").Value = "37.58"
$ws.Range("E14").Value = "  +7.07%  "
$ws.Range("D15").Value = "3.683.29"
$ws.Range("E15").Value = "  +3.71%  "
$ws.Range("D16").Value = "65.278.77"
$ws.Range("E16").Value = "  +2.28%  "
$ws.Range("D17").Value = "3.173.38"
$ws.Range("E17").Value = "  +4.02%  "
$ws.Range("D18").Value = "7.19"
$ws.Range("E18").Value = "  +6.49%  "
$ws.Range("E19").Value = "  +1.25%  "
$ws.Range("D20").Value = "511.62"
$ws.Range("E20").Value = "  +4.68%  "
$ws.Range("D21").Value = "14.88"
$ws.Range("E21").Value = "  +5.20%  "
$ws.Range("D22").Value = "0.726"
$ws.Range("E22").Value = "  +6.53%  "
$ws.Range("D23").Value = "15.32"
$ws.Range("E23").Value = "  +6.23%  "
$ws.Range("D24").Value = "7.83"
$ws.Range("E24").Value = "  +4.28%  "
$ws.Range("D25").Value = "85.12"
$ws.Range("E25").Value = "  +3.17%  "
$ws.Range("D26").Value = "1.00"
$ws.Range("E26").Value = "  +0.13%  "
$ws.Range("D27").Value = "9.09"
$ws.Range("E27").Value = "  +12.39%  "
$ws.Range("D28").Value = "2.94"
$ws.Range("E28").Value = "  +5.20%  "
$ws.Range("D29").Value = "2.20"
$ws.Range("E29").Value = "  +8.37%  "
$ws.Range("D30").Value = "2.84"
$ws.Range("E30").Value = "  +15.72%  "
$ws.Range("D31").Value = "27.90"
$ws.Range("E31").Value = "  +6.25%  "
$ws.Range("D32").Value = "0.999"
$ws.Range("E32").Value = "  -0.03%  "
$ws.Range("E33").Value = "  +4.43%  "
$ws.Range("D34").Value = "6.34"
$ws.Range("E34").Value = "  +12.01%  "
$ws.Range("D35").Value = "6.61"
$ws.Range("E35").Value = "  +6.84%  "
$ws.Range("D36").Value = "55.80"
$ws.Range("E36").Value = "  +1.02%  "
$ws.Range("D37").Value = "0.0903"
$ws.Range("E37").Value = "  +10.95%  "
$ws.Range("D38").Value = "476.44"
$ws.Range("E38").Value = "  +7.97%  "
$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").Value = "0.0421"
$ws.Range("E39").Value = "  +3.23%  "
$ws.Range("B40").Value = "dogwifhat"
$ws.Range("C40").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D40").Value = "3.06"
$ws.Range("E40").Value = "  +10.85%  "
$ws.Range("D41").Value = "8.69"
$ws.Range("E41").Value = "  +4.60%  "
$ws.Range("D42").Value = "3.068.61"
$ws.Range("E42").Value = "  +1.77%  "
$ws.Range("E43").Value = "  +2.87%  "
$ws.Range("D44").Value = "2.48"
$ws.Range("E44").Value = "  +12.00%  "
$ws.Range("D45").Value = "0.286"
$ws.Range("E45").Value = "  +6.06%  "
$ws.Range("D46").Value = "29.24"
$ws.Range("E47").Value = "  +17.95%  "
$ws.Range("E49").Value = "  +1.36%  "
$ws.Range("D50").Value = "2.26"
$ws.Range("E50").Value = "  +8.09%  "
$ws.Range("D51").Value = "120.81"
$ws.Range("E51").Value = "  +2.38%  "
